$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column D: "col"
$ws.Range("D1").Value = "col"

# Copy style from C1 (header style) to D1
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122) # xlPasteFormats

# Data values for new column D
$ws.Range("D2").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("D4").Value = 0

# Update the active cell selection to D5
$ws.Range("D5").Select()
